$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Record the final selection on Sheet1 before moving on (matches the diff:
# Sheet1's selection ends up at B25, tabSelected is dropped once Sheet2 is active).
$ws1.Activate()
$ws1.Range("B25").Select()

# Add the new "Sheet2" after Sheet1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Match Sheet1's page margins (0.75/0.75/1/1/0.5/0.5 in inches -> points).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

$data = @(
    @(1, 1, 1, 1),
    @(2, 2, 2, 2),
    @(3, 3, 3, 3),
    @(4, 4, 4, 4),
    @(5, 5, 5, 5),
    @(0, 11, 22, 33),
    @(7, 7, 7, 7),
    @(78, 78, 78, 78),
    @(89, 91, 95, 100)
)

$r = 3
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws2.Range("B12").Select()
